# BiosampleList.xlsx update
#
# The DATA sheet header row is extended with 7 new leading fields
# (isolate, organism, taxon_id, bio_material, specimen_voucher,
# collected_by, collection date), pushing the existing fields
# (country .. isolation and growth condition) from columns A:J to H:Q.
# The new "collection date" header gets a yyyy-mm-dd date number format.
# The DATA tab becomes the active/selected sheet (it was INFO before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Remember the existing header labels (country .. isolation and growth
# condition) before they get overwritten.
$old = @()
for ($c = 1; $c -le 10; $c++) { $old += $ws.Cells.Item(1, $c).Value() }

# The destination columns K:Q currently have no cells at all, so give
# them the same bold/filled header look as the rest of row 1 (copied
# from the existing header format) before values land on them.
$ws.Range("J1").Copy()
$ws.Range("K1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Shift the existing headers from A1:J1 over to H1:Q1.
for ($c = 1; $c -le 10; $c++) {
  $ws.Cells.Item(1, $c + 7).Value = $old[$c - 1]
}

# Fill in the 7 new headers in A1:G1.
$ws.Cells.Item(1, 1).Value = "isolate"
$ws.Cells.Item(1, 2).Value = "organism"
$ws.Cells.Item(1, 3).Value = "taxon_id"
$ws.Cells.Item(1, 4).Value = "bio_material"
$ws.Cells.Item(1, 5).Value = "specimen_voucher"
$ws.Cells.Item(1, 6).Value = "collected_by"
$ws.Cells.Item(1, 7).Value = "collection date"

# "collection date" is formatted as a date.
$ws.Cells.Item(1, 7).NumberFormat = "yyyy\-mm\-dd;@"

# Widen the columns whose header text changed so the new labels fit.
$ws.Columns("D").ColumnWidth = 18.166666666666668
$ws.Columns("E").ColumnWidth = 15.833333333333334
$ws.Columns("F").ColumnWidth = 15.666666666666666
$ws.Columns("H").ColumnWidth = 18.666666666666668
$ws.Columns("I").ColumnWidth = 21.333333333333332
$ws.Columns("L").ColumnWidth = 10.666666666666666
$ws.Columns("M").ColumnWidth = 10.333333333333334
$ws.Columns("N").ColumnWidth = 19.166666666666668
$ws.Columns("O").ColumnWidth = 28.0
$ws.Columns("P").ColumnWidth = 18.333333333333332
$ws.Columns("Q").ColumnWidth = 29.666666666666668

# DATA becomes the active sheet/tab, with Q6 selected (it was INFO,
# selection C25, before).
$ws.Activate() | Out-Null
$ws.Range("Q6").Select() | Out-Null
